$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells removed entirely (naive forecaster bug fix removed insufficient-data rows)
$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("C4").ClearContents()

# Cells whose forecast values were recalculated after the bug fix
$ws.Range("E3").Value = -0.7976031983999876
$ws.Range("E4").Value = -5.168396053267521
$ws.Range("E5").Value = -8.396348489509165
$ws.Range("C6").Value = -3.956152295564885
$ws.Range("C8").Value = 1.234995474941436
$ws.Range("C10").Value = 0.8993608108207818
$ws.Range("E10").Value = 0.4755443417510108
$ws.Range("C11").Value = 1.15368307467123
$ws.Range("E11").Value = 3.648892256099967
$ws.Range("C13").Value = 0.2186142574756245
$ws.Range("E13").Value = 0.4006004000999486
$ws.Range("C14").Value = 0.02019328874802717
$ws.Range("E14").Value = -1.194807813319176
$ws.Range("E15").Value = -3.55169094390001
$ws.Range("C17").Value = -0.2262139320475476
$ws.Range("E18").Value = -1.02250637024307
$ws.Range("E19").Value = -0.3994003999000184
$ws.Range("E20").Value = -0.7240982069265045
$ws.Range("C21").Value = 0.3239252862367259
$ws.Range("C22").Value = 0.0720185131838802
$ws.Range("E22").Value = 1.255028673974068
$ws.Range("E24").Value = -0.3994003998999962
$ws.Range("C25").Value = 0.347785972938075
$ws.Range("E25").Value = -1.516043567048964
$ws.Range("E26").Value = -3.305525567352907
$ws.Range("C28").Value = -0.4781004700720182
$ws.Range("E28").Value = 0.8024032016000104
$ws.Range("C29").Value = -0.9254001004750156
$ws.Range("E29").Value = -1.738778148048625
$ws.Range("C30").Value = -0.8017595264762423
$ws.Range("E31").Value = -2.378486270399993
$ws.Range("C32").Value = -1.197849743493751
$ws.Range("E32").Value = -3.161804390400014
$ws.Range("E34").Value = 11.45073880931156
$ws.Range("C35").Value = 2.745330323453499
$ws.Range("E35").Value = 2.421686529599998
$ws.Range("C36").Value = 1.098150690304189
$ws.Range("E36").Value = -2.378486270399993
$ws.Range("E37").Value = -2.477884468621794
$ws.Range("C38").Value = 0.9704846793491706
$ws.Range("E38").Value = -1.172596637408219
$ws.Range("E39").Value = 0.8024032015999882
$ws.Range("C40").Value = -2.376072963557374
$ws.Range("E45").Value = -2.011011124776052
$ws.Range("C46").Value = 0.3928252664241683
$ws.Range("C48").Value = 1.052599339874583
$ws.Range("E48").Value = 1.205410808099971
$ws.Range("E49").Value = -0.3246811210723788
$ws.Range("C50").Value = 0.3224026462283369
$ws.Range("E50").Value = -3.037731958703715
$ws.Range("C51").Value = -3.912174452849149
$ws.Range("C52").Value = -0.7109608111999011
$ws.Range("E53").Value = -2.500862062524423
